$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 1.1.1
$ws.Range("B3").Value = "1.1.1"

# Experimental: (empty) -> "false" (stored as literal text, not boolean,
# so prefix with an apostrophe to force text type)
$ws.Range("B7").Value = "'false"

# Date: 2022-05-04T10:16:52-05:00 -> 2022-10-21T09:04:31-05:00
$ws.Range("B8").Value = "2022-10-21T09:04:31-05:00"
